$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 15.571428
$ws.Range("I6").Value = 13.166667
$ws.Range("J6").Value = 30
$ws.Range("K6").Value = 39.500001
$ws.Range("L6").Value = 90
$ws.Range("M6").Value = 72.499999
$ws.Range("N6").Value = -314
$ws.Range("H32").Value = 7058.5
$ws.Range("J32").Value = 6850.6665
$ws.Range("L32").Value = 6850.6665
$ws.Range("N32").Value = -7502.6665
$ws.Range("H55").Value = 419
$ws.Range("I55").Value = 398.33334
$ws.Range("K55").Value = 398.33334
$ws.Range("M55").Value = -184.33334
$ws.Range("H80").Value = 1569.3125
$ws.Range("J80").Value = 1658.091
$ws.Range("L80").Value = 4974.272999999999
$ws.Range("N80").Value = -6970.272999999999
$ws.Range("H83").Value = 1569.3125
$ws.Range("J83").Value = 1658.091
$ws.Range("L83").Value = 14922.819
$ws.Range("N83").Value = -24906.819
$ws.Range("H136").Value = 88749.75
$ws.Range("J136").Value = 88749.75
$ws.Range("L136").Value = 88749.75
$ws.Range("N136").Value = -98949.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H53").Value = 0
$ws.Range("I53").Value = 0
$ws.Range("K53").Value = 0
$ws.Range("M53").ClearContents()
$ws.Range("H61").Value = 13749
$ws.Range("J61").Value = 13749
$ws.Range("L61").Value = 13749
$ws.Range("N61").Value = -14173
$ws.Range("H75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("L75").ClearContents()
$ws.Range("N75").Value = 0
$ws.Range("H78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("L78").ClearContents()
$ws.Range("N78").Value = 0

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 5170
$ws.Range("I20").Value = 3999.3333
$ws.Range("J20").Value = 8682
$ws.Range("K20").Value = 3999.3333
$ws.Range("L20").Value = 8682
$ws.Range("M20").Value = -3752.3333
$ws.Range("N20").Value = -9176
$ws.Range("H22").Value = 184.125
$ws.Range("I22").Value = 184.125
$ws.Range("K22").Value = 184.125
$ws.Range("M22").Value = -11.125
$ws.Range("H54").Value = 4872.5
$ws.Range("I54").Value = 4872.5
$ws.Range("K54").Value = 4872.5
$ws.Range("M54").Value = -4388.5
$ws.Range("H136").Value = 13749
$ws.Range("J136").Value = 13749
$ws.Range("L136").Value = 41247
$ws.Range("N136").Value = -46347

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H47").Value = 18499
$ws.Range("I47").Value = 17998
$ws.Range("J47").Value = 19000
$ws.Range("K47").Value = 17998
$ws.Range("L47").Value = 19000
$ws.Range("M47").Value = -17432
$ws.Range("N47").Value = -20132
$ws.Range("H48").Value = 42400
$ws.Range("J48").Value = 42400
$ws.Range("L48").Value = 42400
$ws.Range("N48").Value = -43352
$ws.Range("H50").Value = 29000
$ws.Range("J50").Value = 29000
$ws.Range("L50").Value = 29000
$ws.Range("N50").Value = -30250
$ws.Range("H51").Value = 6000
$ws.Range("I51").Value = 6000
$ws.Range("K51").Value = 6000
$ws.Range("M51").Value = -5264
$ws.Range("H55").Value = 3280
$ws.Range("I55").Value = 3280
$ws.Range("K55").Value = 3280
$ws.Range("M55").Value = -2965
$ws.Range("H61").Value = 6000
$ws.Range("I61").Value = 6000
$ws.Range("K61").Value = 6000
$ws.Range("M61").Value = -5652
$ws.Range("H68").Value = 45266.668
$ws.Range("I68").Value = 40000
$ws.Range("J68").Value = 47900
$ws.Range("K68").Value = 40000
$ws.Range("L68").Value = 47900
$ws.Range("M68").Value = -39251
$ws.Range("N68").Value = -49398
$ws.Range("H71").Value = 45266.668
$ws.Range("I71").Value = 40000
$ws.Range("J71").Value = 47900
$ws.Range("K71").Value = 120000
$ws.Range("L71").Value = 143700
$ws.Range("M71").Value = -116256
$ws.Range("N71").Value = -151188
$ws.Range("H122").Value = 2306.8
$ws.Range("I122").Value = 2631
$ws.Range("K122").Value = 7893
$ws.Range("M122").Value = -5443

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 741.6667
$ws.Range("I2").Value = 372.86957
$ws.Range("J2").Value = 2862.25
$ws.Range("K2").Value = 2237.21742
$ws.Range("L2").Value = 17173.5
$ws.Range("M2").Value = -2124.21742
$ws.Range("N2").Value = -17399.5
$ws.Range("H34").Value = 458
$ws.Range("I34").Value = 175
$ws.Range("J34").Value = 599.5
$ws.Range("K34").Value = 525
$ws.Range("L34").Value = 1798.5
$ws.Range("M34").Value = -441
$ws.Range("N34").Value = -1966.5
$ws.Range("H39").Value = 2097.5
$ws.Range("I39").Value = 2195
$ws.Range("J39").Value = 2000
$ws.Range("K39").Value = 6585
$ws.Range("L39").Value = 6000
$ws.Range("M39").Value = -6291
$ws.Range("N39").Value = -6588
$ws.Range("H55").Value = 299.5
$ws.Range("I55").Value = 299.5
$ws.Range("K55").Value = 898.5
$ws.Range("M55").Value = -721.5
$ws.Range("H117").Value = 1991.5
$ws.Range("J117").Value = 1718
$ws.Range("L117").Value = 5154
$ws.Range("N117").Value = -12038
$ws.Range("H132").Value = 6500
$ws.Range("I132").Value = 6500
$ws.Range("K132").Value = 19500
$ws.Range("M132").Value = -16970

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 28998
$ws.Range("I15").Value = 28998
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 28998
$ws.Range("L15").ClearContents()
$ws.Range("M15").Value = -28710
$ws.Range("N15").Value = 0
$ws.Range("H81").Value = 28998
$ws.Range("I81").Value = 28998
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 28998
$ws.Range("L81").ClearContents()
$ws.Range("M81").Value = -28000
$ws.Range("N81").Value = 0
$ws.Range("H84").Value = 28998
$ws.Range("I84").Value = 28998
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 86994
$ws.Range("L84").ClearContents()
$ws.Range("M84").Value = -82002
$ws.Range("N84").Value = 0
$ws.Range("H113").Value = 2644.5
$ws.Range("I113").Value = 2644.5
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 2644.5
$ws.Range("L113").Value = 0
$ws.Range("M113").ClearContents()
$ws.Range("N113").Value = -474.5
$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("M126").ClearContents()
$ws.Range("H132").Value = 8200
$ws.Range("I132").Value = 8200
$ws.Range("K132").Value = 24600
$ws.Range("M132").Value = -22070

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2531.2727
$ws.Range("J22").Value = 2285.1428
$ws.Range("L22").Value = 2285.1428
$ws.Range("N22").Value = -2875.1428
$ws.Range("H27").Value = 2531.2727
$ws.Range("J27").Value = 2285.1428
$ws.Range("L27").Value = 2285.1428
$ws.Range("N27").Value = -2499.1428
$ws.Range("H40").Value = 9125
$ws.Range("I40").Value = 9000
$ws.Range("J40").Value = 9250
$ws.Range("K40").Value = 9000
$ws.Range("L40").Value = 9250
$ws.Range("M40").Value = -8864
$ws.Range("N40").Value = -9522
$ws.Range("H55").Value = 1259.6
$ws.Range("J55").Value = 1666.3334
$ws.Range("L55").Value = 1666.3334
$ws.Range("N55").Value = -2012.3334
$ws.Range("H70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("L70").ClearContents()
$ws.Range("N70").Value = 0
$ws.Range("H73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("L73").ClearContents()
$ws.Range("N73").Value = 0
$ws.Range("H100").Value = 2198.75
$ws.Range("I100").Value = 2198.75
$ws.Range("K100").Value = 2198.75
$ws.Range("M100").Value = -1657.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H80").Value = 0
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("L80").ClearContents()
$ws.Range("M80").ClearContents()
$ws.Range("N80").Value = 0
$ws.Range("H83").Value = 0
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("L83").ClearContents()
$ws.Range("M83").ClearContents()
$ws.Range("N83").Value = 0
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").ClearContents()
$ws.Range("N135").Value = 0
$ws.Range("H136").Value = 3500
$ws.Range("I136").Value = 3500
$ws.Range("K136").Value = 10500
$ws.Range("M136").Value = -7950
